$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.759.16"
$ws.Range("E2").Value = "  +1.46%  "

$ws.Range("D3").Value = "2.621.36"
$ws.Range("E3").Value = "  +0.97%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.62"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.82%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.79"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.99%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  +1.69%  "

$ws.Range("D9").Value = "2.619.80"
$ws.Range("E9").Value = "  +0.92%  "

$ws.Range("E10").Value = "  +10.02%  "

$ws.Range("E11").Value = "  +0.83%  "

$ws.Range("E12").Value = "  +0.99%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.13"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.44%  "

$ws.Range("E15").Value = "  +3.62%  "

$ws.Range("D16").Value = "3.099.13"
$ws.Range("E16").Value = "  +1.55%  "

$ws.Range("D17").Value = "67.717.78"
$ws.Range("E17").Value = "  +2.07%  "

$ws.Range("D18").Value = "2.619.52"
$ws.Range("E18").Value = "  +1.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.33"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.28%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "365.69"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.71%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.63"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.51%  "

$ws.Range("E22").Value = "  -0.43%  "

$ws.Range("E23").Value = "  +4.21%  "

$ws.Range("E24").Value = "  +0.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.16"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +4.23%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.13"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.36%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000105"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.02%  "

$ws.Range("D28").Value = "2.744.87"
$ws.Range("E28").Value = "  +0.61%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "586.07"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.91%  "

$ws.Range("E30").Value = "  +0.34%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.44"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.66%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.96"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.80%  "

$ws.Range("E33").Value = "  +0.47%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.131"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.90%  "

$ws.Range("E35").Value = "  +0.07%  "

$ws.Range("E36").Value = "  -1.33%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.98"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.43%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.47"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.95%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "155.60"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.06%  "

$ws.Range("E40").Value = "  +1.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.43"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.27%  "

$ws.Range("E42").Value = "  +3.16%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.65"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.66%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.13"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.45%  "

$ws.Range("E45").Value = "  -0.10%  "

$ws.Range("E46").Value = "  +0.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "157.38"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.44%  "

$ws.Range("D48").Value = "0.0₆0288"
$ws.Range("E48").Value = "  -7.03%  "

$ws.Range("E49").Value = "  +0.20%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "21.10"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.60%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.625"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.13%  "
